$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "93.987.90"
$ws.Range("E2").Value = "  +1.12%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.074.60"
$ws.Range("E3").Value = "  -1.44%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.97"
$ws.Range("E5").Value = "  -3.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "607.93"
$ws.Range("E6").Value = "  -1.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.09"
$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.378"
$ws.Range("E8").Value = "  -5.68%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.821"
$ws.Range("E10").Value = "  +12.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.072.28"
$ws.Range("E11").Value = "  -1.37%  "

$ws.Range("E12").Value = "  -3.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "93.870.27"
$ws.Range("E13").Value = "  +1.43%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000239"
$ws.Range("E14").Value = "  -5.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.91"
$ws.Range("E15").Value = "  -1.69%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.26"
$ws.Range("E16").Value = "  -4.47%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.645.15"
$ws.Range("E17").Value = "  -1.36%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.081.35"
$ws.Range("E18").Value = "  -0.84%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.64"
$ws.Range("E19").Value = "  -1.91%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.50"
$ws.Range("E20").Value = "  -1.62%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.73"
$ws.Range("E21").Value = "  -1.53%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "440.03"
$ws.Range("E22").Value = "  -1.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.80"
$ws.Range("E23").Value = "  -6.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000191"
$ws.Range("E24").Value = "  -6.53%  "

$ws.Range("E25").Value = "  +5.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "5.52"
$ws.Range("E26").Value = "  -4.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "84.86"
$ws.Range("E27").Value = "  -2.36%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.89"
$ws.Range("E28").Value = "  +1.09%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.252.00"
$ws.Range("E29").Value = "  -0.86%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.248"
$ws.Range("E31").Value = "  +7.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.177"
$ws.Range("E32").Value = "  +5.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.122"
$ws.Range("E33").Value = "  -10.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.12"
$ws.Range("E34").Value = "  -2.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.994"
$ws.Range("E35").Value = "  -0.56%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.73"
$ws.Range("E36").Value = "  -4.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.157"
$ws.Range("E37").Value = "  -4.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.58"
$ws.Range("E38").Value = "  -2.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.87"
$ws.Range("E39").Value = "  -2.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.443"
$ws.Range("E40").Value = "  +0.92%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "23.96"
$ws.Range("E41").Value = "  +3.66%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.26"
$ws.Range("E42").Value = "  -3.75%  "

$ws.Range("B43").Value = "MantraDAO"
$ws.Range("C43").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.69"
$ws.Range("E43").Value = "  -13.34%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "463.95"
$ws.Range("E44").Value = "  -4.23%  "

$ws.Range("E45").Value = "  +0.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.10"
$ws.Range("E46").Value = "  -11.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "160.64"
$ws.Range("E47").Value = "  -0.15%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.84"
$ws.Range("E48").Value = "  -4.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.674"
$ws.Range("E49").Value = "  -3.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.73"
$ws.Range("E50").Value = "  -0.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.997"
$ws.Range("E51").Value = "  +0.09%  "
